$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 794
$ws1.Range("F5").Value = 156
$ws1.Range("F6").Value = 20
$ws1.Range("F7").Value = 178
$ws1.Range("F8").Value = 361
$ws1.Range("F9").Value = 475
$ws1.Range("F12").Value = 12083
$ws1.Range("F13").Value = 5443

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 794
$ws4.Range("F7").Value = 156
$ws4.Range("F8").Value = 20
$ws4.Range("F9").Value = 178
$ws4.Range("F10").Value = 361
$ws4.Range("F11").Value = 475
$ws4.Range("F14").Value = 12083
$ws4.Range("F16").Value = 5443

$wb.Save()
